$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revert C15 and C16 status from "done" back to "open"
$ws.Range("C15").Value = "open"
$ws.Range("C16").Value = "open"

# Revert row 19 (item 17, "Import & Export reg setting" details) back to an
# empty/open placeholder row, matching the formatting of the rows around it.
$ws.Range("B20:E20").Copy()
$ws.Range("B19:E19").PasteSpecial(-4122)
$ws.Range("B19").ClearContents()
$ws.Range("C19").Value = "open"
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Rows.Item(19).EntireRow.AutoFit()

# Selection moves to B21
$ws.Range("B21").Select()

# Column D's custom width is no longer needed; column E narrows slightly.
$ws.Columns.Item(4).ColumnWidth = 8.43
$ws.Columns.Item(5).ColumnWidth = 49.7
